$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 7 and 8, shifting the existing rows 7..59 down to 9..61
$ws.Rows("7:8").Insert()

# New row 7: Choclo / Choclero / Primera, Vega Monumental Concepción, 2021-12-22
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 44552
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = 100112024
$ws.Cells.Item(7, 7).Value = "Choclo"
$ws.Cells.Item(7, 8).Value = "Choclero"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 22000
$ws.Cells.Item(7, 12).Value = 24000
$ws.Cells.Item(7, 13).Value = 23000
$ws.Cells.Item(7, 14).Value = "$/malla 50 unidades"
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 16).Value = 460
$ws.Cells.Item(7, 17).Value = 50
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# New row 8: Choclo / Choclero / Primera, Vega Monumental Concepción, 2021-12-22
$ws.Cells.Item(8, 1).Value = 11
$ws.Cells.Item(8, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(8, 3).Value = "Bíobío"
$ws.Cells.Item(8, 4).Value = 44552
$ws.Cells.Item(8, 5).Value = 8
$ws.Cells.Item(8, 6).Value = 100112024
$ws.Cells.Item(8, 7).Value = "Choclo"
$ws.Cells.Item(8, 8).Value = "Choclero"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 11000
$ws.Cells.Item(8, 12).Value = 12000
$ws.Cells.Item(8, 13).Value = 11500
$ws.Cells.Item(8, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(8, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value = 164
$ws.Cells.Item(8, 17).Value = 70
$ws.Cells.Item(8, 18).Value = "Hortaliza"
